# Fruta / hortaliza, semanal
# Insert a new weekly record at row 4 ("Macroferia Regional de Talca" - Chirimoya),
# pushing the existing rows 4-10 down to rows 5-11, and populate the new row 4
# with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 (inherits formatting, e.g. the date
# number format on column D, from the row above - matching the rest of the sheet).
$ws.Rows(4).Insert()

$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(4, 3).Value = "Maule"
$ws.Cells.Item(4, 4).Value = 44447
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100107
$ws.Cells.Item(4, 8).Value = "Otros"
$ws.Cells.Item(4, 9).Value = 100107002
$ws.Cells.Item(4, 10).Value = "Chirimoya"
$ws.Cells.Item(4, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(4, 12).Value = "Especial"
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 32000
$ws.Cells.Item(4, 15).Value = 32000
$ws.Cells.Item(4, 16).Value = 32000
$ws.Cells.Item(4, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 3200
$ws.Cells.Item(4, 20).Value = 10
